$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 383.9
$ws.Range("I33").Value = 393.93332
$ws.Range("J33").Value = 353.8
$ws.Range("K33").Value = 393.93332
$ws.Range("L33").Value = 353.8
$ws.Range("M33").Value = -164.93332
$ws.Range("N33").Value = -811.8
# Row 41
$ws.Range("H41").Value = 1579.25
$ws.Range("I41").Value = 1156.5555
$ws.Range("J41").Value = 2847.3333
$ws.Range("K41").Value = 1156.5555
$ws.Range("L41").Value = 2847.3333
$ws.Range("M41").Value = -716.5554999999999
$ws.Range("N41").Value = -3727.3333
# Row 52
$ws.Range("H52").Value = 8083.1665
$ws.Range("I52").Value = 8083.1665
$ws.Range("K52").Value = 24249.4995
$ws.Range("M52").Value = -24089.4995
# Row 53
$ws.Range("H53").Value = 542.5833
$ws.Range("I53").Value = 209.8
$ws.Range("K53").Value = 209.8
$ws.Range("M53").Value = 427.2
# Row 100
$ws.Range("H100").Value = 6467.4
$ws.Range("I100").Value = 5035
$ws.Range("J100").Value = 7899.8
$ws.Range("K100").Value = 5035
$ws.Range("L100").Value = 7899.8
$ws.Range("M100").Value = -4494
$ws.Range("N100").Value = -8981.799999999999
# Row 133
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 20005634
$ws.Range("I32").Value = 20839034
$ws.Range("K32").Value = 20839034
$ws.Range("M32").Value = -20838747
# Row 45
$ws.Range("H45").Value = 2762.55
$ws.Range("I45").Value = 2739.7058
$ws.Range("K45").Value = 2739.7058
$ws.Range("M45").Value = -2362.7058
# Row 61
$ws.Range("H61").Value = 4706.857
$ws.Range("I61").Value = 2999.6667
$ws.Range("K61").Value = 2999.6667
$ws.Range("M61").Value = -2787.6667
# Row 102
$ws.Range("H102").Value = 9524923
$ws.Range("I102").Value = 1167.25
$ws.Range("K102").Value = 1167.25
$ws.Range("M102").Value = 454.75
# Row 110
$ws.Range("H110").Value = 933.5833
$ws.Range("I110").Value = 984.65
$ws.Range("K110").Value = 984.65
$ws.Range("M110").Value = 1060.35
# Row 122
$ws.Range("H122").Value = 2156.7827
$ws.Range("I122").Value = 1976.8235
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 5930.470499999999
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -3480.470499999999
$ws.Range("N122").Value = -12900.0001
# Row 136
$ws.Range("H136").Value = 4706.857
$ws.Range("I136").Value = 2999.6667
$ws.Range("K136").Value = 8999.000100000001
$ws.Range("M136").Value = -6449.000100000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 5813.3335
$ws.Range("I105").Value = 3586.6667
$ws.Range("K105").Value = 3586.6667
$ws.Range("M105").Value = -1839.6667
# Row 107
$ws.Range("H107").Value = 2230.3333
$ws.Range("I107").Value = 2230.3333
$ws.Range("K107").Value = 2230.3333
$ws.Range("M107").Value = -310.3332999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 971.3889
$ws.Range("I22").Value = 370.8
$ws.Range("K22").Value = 370.8
$ws.Range("M22").Value = -20.80000000000001
# Row 31
$ws.Range("H31").Value = 1741.2222
$ws.Range("I31").Value = 1728.6857
$ws.Range("K31").Value = 1728.6857
$ws.Range("M31").Value = -1433.6857
# Row 34
$ws.Range("H34").Value = 1741.2222
$ws.Range("I34").Value = 1728.6857
$ws.Range("K34").Value = 1728.6857
$ws.Range("M34").Value = -1526.6857
# Row 35
$ws.Range("H35").Value = 2137.5
$ws.Range("I35").Value = 1850
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 1850
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = -1556
$ws.Range("N35").Value = -3588
# Row 58
$ws.Range("H58").Value = 3244.9375
$ws.Range("I58").Value = 2755
$ws.Range("K58").Value = 2755
$ws.Range("M58").Value = -2552
# Row 86
$ws.Range("H86").Value = 139906
$ws.Range("I86").Value = 299999
$ws.Range("J86").Value = 59859.5
$ws.Range("K86").Value = 299999
$ws.Range("L86").Value = 59859.5
$ws.Range("M86").Value = -298876
$ws.Range("N86").Value = -62105.5
# Row 89
$ws.Range("H89").Value = 139906
$ws.Range("I89").Value = 299999
$ws.Range("J89").Value = 59859.5
$ws.Range("K89").Value = 1499995
$ws.Range("L89").Value = 299297.5
$ws.Range("M89").Value = -1494379
$ws.Range("N89").Value = -310529.5
# Row 99
$ws.Range("H99").Value = 28056342
$ws.Range("I99").Value = 4070120.8
$ws.Range("K99").Value = 4070120.8
$ws.Range("M99").Value = -4068622.8
# Row 105
$ws.Range("H105").Value = 19199.8
$ws.Range("I105").Value = 3666.3333
$ws.Range("K105").Value = 3666.3333
$ws.Range("M105").Value = -1919.3333
# Row 122
$ws.Range("H122").Value = 734890.5600000001
$ws.Range("I122").Value = 2043293.6
$ws.Range("K122").Value = 6129880.800000001
$ws.Range("M122").Value = -6127430.800000001
# Row 126
$ws.Range("H126").Value = 28056342
$ws.Range("I126").Value = 4070120.8
$ws.Range("K126").Value = 12210362.4
$ws.Range("M126").Value = -12207892.4
# Row 132
$ws.Range("H132").Value = 4217.5
$ws.Range("I132").Value = 4217.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12652.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10122.5
$ws.Range("N132").ClearContents()
# Row 135
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140
# Row 136
$ws.Range("H136").Value = 3244.9375
$ws.Range("I136").Value = 2755
$ws.Range("K136").Value = 8265
$ws.Range("M136").Value = -5715

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 17964806
$ws.Range("I4").Value = 25150538
$ws.Range("K4").Value = 75451614
$ws.Range("M4").Value = -75451502
# Row 41
$ws.Range("H41").Value = 534.61536
$ws.Range("J41").Value = 155
$ws.Range("L41").Value = 465
$ws.Range("N41").Value = -1141
# Row 63
$ws.Range("H63").Value = 4832.385
$ws.Range("I63").Value = 3705.5
$ws.Range("K63").Value = 11116.5
$ws.Range("M63").Value = -10367.5
# Row 64
$ws.Range("H64").Value = 2874.5
# Row 66
$ws.Range("H66").Value = 4832.385
$ws.Range("I66").Value = 3705.5
$ws.Range("K66").Value = 33349.5
$ws.Range("M66").Value = -29605.5
# Row 67
$ws.Range("H67").Value = 2874.5
# Row 70
$ws.Range("H70").Value = 449448.5
$ws.Range("I70").Value = 449448
$ws.Range("K70").Value = 1348344
$ws.Range("M70").Value = -1348029
# Row 73
$ws.Range("H73").Value = 449448.5
$ws.Range("I73").Value = 449448
$ws.Range("K73").Value = 1348344
$ws.Range("M73").Value = -1347252
# Row 92
$ws.Range("H92").Value = 711.1111
$ws.Range("I92").Value = 483.33334
$ws.Range("J92").Value = 1166.6666
$ws.Range("K92").Value = 1450.00002
$ws.Range("L92").Value = 3499.9998
$ws.Range("M92").Value = -202.0000199999999
$ws.Range("N92").Value = -5995.9998
# Row 98
$ws.Range("H98").Value = 1599.2
$ws.Range("I98").Value = 1599.2
$ws.Range("K98").Value = 4797.6
$ws.Range("M98").Value = -3299.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 1000.7692
$ws.Range("I107").Value = 1371.375
$ws.Range("J107").Value = 407.8
$ws.Range("K107").Value = 1371.375
$ws.Range("L107").Value = 407.8
$ws.Range("M107").Value = 548.625
$ws.Range("N107").Value = -4247.8
# Row 133
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 24
$ws.Range("H24").Value = 7602.5
$ws.Range("I24").Value = 7602.5
$ws.Range("K24").Value = 7602.5
$ws.Range("M24").Value = -7259.5
# Row 25
$ws.Range("H25").Value = 1050
$ws.Range("J25").Value = 1050
$ws.Range("L25").Value = 1050
$ws.Range("N25").Value = -1510
# Row 132
$ws.Range("H132").Value = 4016.5908
$ws.Range("I132").Value = 2799.1892
$ws.Range("J132").Value = 10451.429
$ws.Range("K132").Value = 8397.567599999998
$ws.Range("L132").Value = 31354.287
$ws.Range("M132").Value = -5867.567599999998
$ws.Range("N132").Value = -36414.287
# Row 136
$ws.Range("H136").Value = 4766.2856
$ws.Range("I136").Value = 3538.2104
$ws.Range("J136").Value = 7358.8887
$ws.Range("K136").Value = 10614.6312
$ws.Range("L136").Value = 22076.6661
$ws.Range("M136").Value = -8064.6312
$ws.Range("N136").Value = -27176.6661

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 15628294
$ws.Range("I62").Value = 2475
$ws.Range("J62").Value = 17860554
$ws.Range("K62").Value = 2475
$ws.Range("L62").Value = 17860554
$ws.Range("M62").Value = -1851
$ws.Range("N62").Value = -17861802
# Row 65
$ws.Range("H65").Value = 15628294
$ws.Range("I65").Value = 2475
$ws.Range("J65").Value = 17860554
$ws.Range("K65").Value = 12375
$ws.Range("L65").Value = 89302770
$ws.Range("M65").Value = -9255
$ws.Range("N65").Value = -89309010
# Row 96
$ws.Range("H96").Value = 66137.25
$ws.Range("I96").Value = 103200.6
$ws.Range("J96").Value = 4365
$ws.Range("K96").Value = 103200.6
$ws.Range("L96").Value = 4365
$ws.Range("M96").Value = -101827.6
$ws.Range("N96").Value = -7111
# Row 100
$ws.Range("H100").Value = 595.6316
$ws.Range("I100").Value = 564.5333000000001
$ws.Range("K100").Value = 1129.0666
$ws.Range("M100").Value = -588.0666000000001
# Row 132
$ws.Range("H132").Value = 2392.739
$ws.Range("I132").Value = 2392.739
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7178.217000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4648.217000000001
$ws.Range("N132").ClearContents()
# Row 136
$ws.Range("H136").Value = 900.7320999999999
$ws.Range("I136").Value = 481.86047
$ws.Range("K136").Value = 1445.58141
$ws.Range("M136").Value = 1104.41859

Write-Host "Applied all edits"